$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 17.05443219947866
$ws.Range("A2").Value = 17.11297689133074
$ws.Range("A3").Value = 16.78313955193191
$ws.Range("A4").Value = 16.33713498461266
$ws.Range("A5").Value = 17.30922597452606
$ws.Range("A6").Value = 16.68082997295395
$ws.Range("A7").Value = 16.32995891797985
$ws.Range("A8").Value = 17.10101204440819
$ws.Range("A9").Value = 16.80536188956136
$ws.Range("A10").Value = 17.32481890889563
